$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This metadata sheet describes, per column (A..J), four facts stacked in
# rows below the human-readable header (row 1):
#   row 2 -> machine slug            (NEW)
#   row 3 -> measure/dimension id    (was row 2)
#   row 4 -> "medida"/"dim" marker   (was row 3)
#   row 5 -> xsd datatype            (was row 4, now fully populated)
#
# A new slug row is introduced right under the headers so that two
# columns can later be related to build hierarchical SKOS concepts. This
# pushes the previous rows 2-4 down one position, and the old, mostly
# empty last row (which only held "mapping-ano.xlsx" in column H) is
# replaced outright by a fully populated xsd-types row.

# Row 2 (new): machine-friendly slug derived from each column's header.
$ws.Range("A2").Value = "edad-grandes-grupos"
$ws.Range("B2").Value = "personas"
$ws.Range("C2").Value = "residencia-provincia-codigo"
$ws.Range("D2").Value = "residencia-comarca-codigo"
$ws.Range("E2").Value = "residencia-ccaa-nombre"
$ws.Range("F2").Value = "residencia-comarca-nombre"
$ws.Range("G2").Value = "residencia-provincia-nombre"
$ws.Range("H2").Value = "ano"
$ws.Range("I2").Value = "sexo"
$ws.Range("J2").Value = "relacion-lugar-de-residencia-y-nacimiento"

# Row 3 (was row 2): measure/dimension identifiers.
$ws.Range("A3").Value = "iaest-measure:edad-grandes-grupos"
$ws.Range("B3").Value = "iaest-measure:personas"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("F3").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("G3").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("H3").Value = "sdmx-dimension:refPeriod"
$ws.Range("I3").Value = "iaest-measure:sexo"
$ws.Range("J3").Value = "iaest-measure:relacion-lugar-de-residencia-y-nacimiento"

# Row 4 (was row 3): medida/dim markers.
$ws.Range("A4").Value = "medida"
$ws.Range("B4").Value = "medida"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "medida"
$ws.Range("F4").Value = "medida"
$ws.Range("G4").Value = "medida"
$ws.Range("H4").Value = "dim"
$ws.Range("I4").Value = "medida"
$ws.Range("J4").Value = "medida"

# Row 5 (was row 4's xsd types; column H previously held the lone value
# "mapping-ano.xlsx" - the row is now fully populated with datatypes).
$ws.Range("A5").Value = "xsd:string"
$ws.Range("B5").Value = "xsd:int"
$ws.Range("C5").Value = "null"
$ws.Range("D5").Value = "null"
$ws.Range("E5").Value = "xsd:string"
$ws.Range("F5").Value = "xsd:string"
$ws.Range("G5").Value = "xsd:string"
$ws.Range("H5").Value = "xsd:date"
$ws.Range("I5").Value = "xsd:string"
$ws.Range("J5").Value = "xsd:string"

# Match the existing "plain" cell formatting used throughout the sheet
# (style index 1) for every newly-populated cell, using column H's
# already-correctly-styled row 5 cell as the format donor.
$ws.Range("H5").Copy()
$ws.Range("A2:G5").PasteSpecial(-4122)
$ws.Range("I5:J5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
